$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @(1, 19643, 5281, 9069, 0, 5814, 0, 31668, 19387),
  @(2, 19040, 5040, 8574, 0, 5639, 0, 15262, 18432),
  @(3, 18162, 4815, 8241, 0, 5394, 0, 14814, 17675),
  @(4, 17939, 4600, 8067, 0, 5206, 0, 14538, 16824),
  @(5, 17043, 4533, 7955, 0, 5166, 0, 16272, 16740),
  @(6, 17442, 4661, 8396, 0, 5324, 0, 26896, 16961),
  @(7, 17016, 4457, 8400, 0, 6436, 0, 13134, 16774),
  @(8, 18140, 4331, 9655, 0, 8429, 0, 14975, 18225),
  @(9, 21348, 4887, 11339, 0, 10065, 0, 15898, 19443),
  @(10, 23261, 5220, 16300, 0, 12623, 0, 16704, 21165),
  @(11, 24338, 5393, 16277, 0, 12894, 0, 17187, 21851),
  @(12, 23109, 5507, 16046, 0, 12282, 0, 16479, 21592),
  @(13, 24018, 5301, 16739, 0, 17565, 0, 17534, 22061),
  @(14, 24893, 5711, 16754, 0, 13142, 0, 17148, 23071),
  @(15, 25292, 5665, 16075, 0, 13332, 0, 16830, 22368),
  @(16, 24339, 5402, 16202, 0, 13289, 0, 16510, 22358),
  @(17, 23612, 5081, 15610, 0, 12948, 0, 16003, 21658),
  @(18, 23281, 4908, 15366, 0, 12435, 0, 15753, 21468),
  @(19, 23789, 5575, 15309, 0, 11345, 0, 16857, 22267),
  @(20, 22074, 5636, 15093, 0, 10183, 0, 16352, 21896),
  @(21, 21667, 5851, 14538, 0, 8028, 0, 16625, 22048),
  @(22, 21596, 5973, 12506, 0, 7058, 0, 17049, 20872),
  @(23, 21121, 5844, 9627, 0, 6232, 0, 16818, 20534),
  @(24, 20491, 5652, 8859, 0, 5806, 0, 16249, 19728)
)

for ($i = 0; $i -lt $data.Length; $i++) {
  $row = $i + 2
  $vals = $data[$i]
  for ($c = 0; $c -lt $vals.Length; $c++) {
    $col = $c + 1
    $ws.Cells.Item($row, $col).Value = $vals[$c]
  }
}
